$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "('Monstrous Hound', ['{3}{R}', 'Creature " + [char]0x2014 + " Dog', 'Monstrous Hound can" + [char]0x2019 + "t attack unless you control more lands than defending player.', 'Monstrous Hound can" + [char]0x2019 + "t block unless you control more lands than attacking player.', '4/4'])"

$ws.Range("A3:A7").ClearContents() | Out-Null
$ws.Range("A2").Value = $newValue
